# Registree stats backup on Wed 28 Apr 2021 13:30:17 SAST
#
# - Two new registrees ("du Plooy", Tammy & Marc, Alberton, 410E) were added
#   to the MD410 Attendance sheet, sorted alphabetically right before the
#   existing "du Toit" row (row 223), pushing every following row down by 2.
# - The "Number of attendees" / "Number of voters" summary rows at the
#   bottom were bumped accordingly (241 -> 243 attendees, 96 -> 98 voters).
# - The "as of <date>" timestamp in row 1 of every sheet was refreshed.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)   # "MD410 Attendance"

# Insert two blank rows right above the current row 223 ("du Toit") so the
# new "du Plooy" entries keep the list in alphabetical order.
$ws.Range("A223:A224").EntireRow.Insert()

# Match the look of the surrounding data rows (25pt tall, thin box border,
# same style index "3" used by every other registree row).
$newRows = $ws.Range("A223:F224")
$newRows.RowHeight = 25
$newRows.Borders.LineStyle = 1

$ws.Cells.Item(223, 1).Value = "du Plooy"
$ws.Cells.Item(223, 2).Value = "Tammy"
$ws.Cells.Item(223, 3).Value = "Alberton"
$ws.Cells.Item(223, 4).Value = "No"
$ws.Cells.Item(223, 5).Value = "Yes"
$ws.Cells.Item(223, 6).Value = "410E"

$ws.Cells.Item(224, 1).Value = "du Plooy"
$ws.Cells.Item(224, 2).Value = "Marc"
$ws.Cells.Item(224, 3).Value = "Alberton"
$ws.Cells.Item(224, 4).Value = "No"
$ws.Cells.Item(224, 5).Value = "Yes"
$ws.Cells.Item(224, 6).Value = "410E"

# The summary rows shifted from 244/245 to 246/247; refresh their counts.
$ws.Cells.Item(246, 1).Value = "Number of attendees: 243"
$ws.Cells.Item(247, 1).Value = "Number of voters: 98"

# Refresh the "as of" timestamp banner on every sheet.
$ws.Cells.Item(1, 1).Value = "MD410 Registrees as of 28/04/2021 13:30"

$ws2 = $wb.Worksheets.Item(2)  # "410E Attendance"
$ws2.Cells.Item(1, 1).Value = "410E Registrees as of 28/04/2021 13:30"

$ws3 = $wb.Worksheets.Item(3)  # "410W Attendance"
$ws3.Cells.Item(1, 1).Value = "410W Registrees as of 28/04/2021 13:30"

$ws4 = $wb.Worksheets.Item(4)  # "410E Voting"
$ws4.Cells.Item(1, 1).Value = "410E Voting details as of 28/04/2021 13:30"

$ws5 = $wb.Worksheets.Item(5)  # "410W Voting"
$ws5.Cells.Item(1, 1).Value = "410W Voting details as of 28/04/2021 13:30"
